# Mexico Liga MX workbook update
# - Several pairs of adjacent data rows get their full record (columns B:AB)
#   swapped with one another (column A, the sequential row index, stays put).
# - The last row (339), representing a not-yet-played fixture, gets updated
#   with a new match id, new date, home/away teams swapped, and refreshed
#   opening/closing odds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $row1, $row2, $firstCol, $lastCol) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $c1 = $ws.Cells.Item($row1, $col)
        $c2 = $ws.Cells.Item($row2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value2 = $v2
        $c2.Value2 = $v1
    }
}

# Row pairs whose full records (columns B..AB) are swapped.
$pairs = @(
    @(34, 35),
    @(72, 73),
    @(98, 99),
    @(128, 129),
    @(130, 131),
    @(200, 201),
    @(222, 223),
    @(237, 238),
    @(251, 252),
    @(282, 283),
    @(303, 304),
    @(318, 319),
    @(322, 323)
)

foreach ($pair in $pairs) {
    Swap-Rows $ws $pair[0] $pair[1] 2 28   # columns B (2) through AB (28)
}

# Row 339: new upcoming fixture data.
# Match id (column B) is stored as text in this workbook, so force a text
# number format before assigning, then clear the format override so the
# cell keeps the same (default) style it had before.
$idCell = $ws.Range("B339")
$idCell.NumberFormat = "@"
$idCell.Value2 = "8241446"
$idCell.ClearFormats()

$ws.Range("D339").Value2 = 45438.91666666666   # 2024-05-26 22:00:00

$ws.Range("E339").Value2 = "Club America"
$ws.Range("F339").Value2 = "Cruz Azul"

$ws.Range("J339").Value2 = 1.8
$ws.Range("K339").Value2 = 3.5
$ws.Range("L339").Value2 = 4
$ws.Range("M339").Value2 = 1.95
$ws.Range("N339").Value2 = 3.3
$ws.Range("O339").Value2 = 4
$ws.Range("P339").Value2 = -0.5
$ws.Range("Q339").Value2 = 1.975
$ws.Range("R339").Value2 = 1.875
$ws.Range("S339").Value2 = 2.25
$ws.Range("T339").Value2 = 2.05
$ws.Range("U339").Value2 = 1.8
